$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Remove the two rows that drop out of the dataset (118e24 at row 6, 119e24 at row 12).
# Delete in descending row order so earlier row indices stay valid.
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(6).Delete()

# Update remaining rows 1-10 with the refreshed computed values (columns B,D,E,F,G,H,I,J; A and C are unchanged)
$ws.Cells.Item(1, 2).Value = 2226
$ws.Cells.Item(1, 4).Value = 1404
$ws.Cells.Item(1, 5).Value = 821
$ws.Cells.Item(1, 6).Value = 873
$ws.Cells.Item(1, 7).Value = 61.6600790513834
$ws.Cells.Item(1, 8).Value = 63.10112359550562
$ws.Cells.Item(1, 9).Value = 0.7436347673397717
$ws.Cells.Item(1, 10).Value = 83.09997510910034

$ws.Cells.Item(2, 2).Value = 2765
$ws.Cells.Item(2, 4).Value = 2072
$ws.Cells.Item(2, 5).Value = 692
$ws.Cells.Item(2, 6).Value = 205
$ws.Cells.Item(2, 7).Value = 90.99692577953448
$ws.Cells.Item(2, 8).Value = 74.96382054992765
$ws.Cells.Item(2, 9).Value = 0.393766461808604
$ws.Cells.Item(2, 10).Value = 75.49118638038635

$ws.Cells.Item(3, 2).Value = 2634
$ws.Cells.Item(3, 4).Value = 2213
$ws.Cells.Item(3, 5).Value = 420
$ws.Cells.Item(3, 6).Value = 64
$ws.Cells.Item(3, 7).Value = 97.18928414580589
$ws.Cells.Item(3, 8).Value = 84.04861374857578
$ws.Cells.Item(3, 9).Value = 0.2124670763827919
$ws.Cells.Item(3, 10).Value = 77.21280813217163

$ws.Cells.Item(4, 2).Value = 2466
$ws.Cells.Item(4, 4).Value = 2277
$ws.Cells.Item(4, 5).Value = 188
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 100
$ws.Cells.Item(4, 8).Value = 92.37322515212982
$ws.Cells.Item(4, 9).Value = 0.08252853380158033
$ws.Cells.Item(4, 10).Value = 78.55643463134766

$ws.Cells.Item(5, 2).Value = 2331
$ws.Cells.Item(5, 4).Value = 2277
$ws.Cells.Item(5, 5).Value = 53
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 100
$ws.Cells.Item(5, 8).Value = 97.72532188841201
$ws.Cells.Item(5, 9).Value = 0.02326602282704126
$ws.Cells.Item(5, 10).Value = 78.6871497631073

$ws.Cells.Item(6, 2).Value = 2499
$ws.Cells.Item(6, 4).Value = 1591
$ws.Cells.Item(6, 5).Value = 906
$ws.Cells.Item(6, 6).Value = 395
$ws.Cells.Item(6, 7).Value = 80.11077542799597
$ws.Cells.Item(6, 8).Value = 63.71645975170204
$ws.Cells.Item(6, 9).Value = 0.6547559134373427
$ws.Cells.Item(6, 10).Value = 87.67913198471069

$ws.Cells.Item(7, 2).Value = 2538
$ws.Cells.Item(7, 4).Value = 1851
$ws.Cells.Item(7, 5).Value = 685
$ws.Cells.Item(7, 6).Value = 135
$ws.Cells.Item(7, 7).Value = 93.20241691842901
$ws.Cells.Item(7, 8).Value = 72.98895899053628
$ws.Cells.Item(7, 9).Value = 0.4126824358329139
$ws.Cells.Item(7, 10).Value = 83.81719088554382

$ws.Cells.Item(8, 2).Value = 2388
$ws.Cells.Item(8, 4).Value = 1967
$ws.Cells.Item(8, 5).Value = 419
$ws.Cells.Item(8, 6).Value = 19
$ws.Cells.Item(8, 7).Value = 99.04330312185297
$ws.Cells.Item(8, 8).Value = 82.43922883487008
$ws.Cells.Item(8, 9).Value = 0.2204328132863614
$ws.Cells.Item(8, 10).Value = 86.12558889389038

$ws.Cells.Item(9, 2).Value = 2243
$ws.Cells.Item(9, 4).Value = 1986
$ws.Cells.Item(9, 5).Value = 255
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 100
$ws.Cells.Item(9, 8).Value = 88.62115127175367
$ws.Cells.Item(9, 9).Value = 0.128334172118772
$ws.Cells.Item(9, 10).Value = 89.00701522827148

$ws.Cells.Item(10, 2).Value = 2057
$ws.Cells.Item(10, 4).Value = 1986
$ws.Cells.Item(10, 5).Value = 69
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 100
$ws.Cells.Item(10, 8).Value = 96.64233576642336
$ws.Cells.Item(10, 9).Value = 0.03472571716155008
$ws.Cells.Item(10, 10).Value = 93.05727291107178
